$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records need to be inserted just before the current
# row 573, pushing the existing rows (old 573-613) down to 575-615.
$ws.Rows("573:574").Insert()

# New row 573: Kiwi Hayward "Especial" quality record for 2021-11-16
$ws.Range("A573").Value = 6
$ws.Range("B573").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C573").Value = "Metropolitana"
$ws.Range("D573").Value = 44516
$ws.Range("E573").Value = 13
$ws.Range("F573").Value = "Fruta"
$ws.Range("G573").Value = 100101
$ws.Range("H573").Value = "Berries"
$ws.Range("I573").Value = 100101007
$ws.Range("J573").Value = "Kiwi"
$ws.Range("K573").Value = "Hayward"
$ws.Range("L573").Value = "Especial"
$ws.Range("M573").Value = 21
$ws.Range("N573").Value = 430000
$ws.Range("O573").Value = 460000
$ws.Range("P573").Value = 441429
$ws.Range("Q573").Value = "`$/bins (450 kilos)"
$ws.Range("R573").Value = "Región de O'Higgins"
$ws.Range("S573").Value = 981
$ws.Range("T573").Value = 450

# New row 574: Kiwi Hayward "Extra (doble especial)" quality record for 2021-03-15
$ws.Range("A574").Value = 6
$ws.Range("B574").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C574").Value = "Metropolitana"
$ws.Range("D574").Value = 44270
$ws.Range("E574").Value = 13
$ws.Range("F574").Value = "Fruta"
$ws.Range("G574").Value = 100101
$ws.Range("H574").Value = "Berries"
$ws.Range("I574").Value = 100101007
$ws.Range("J574").Value = "Kiwi"
$ws.Range("K574").Value = "Hayward"
$ws.Range("L574").Value = "Extra (doble especial)"
$ws.Range("M574").Value = 10
$ws.Range("N574").Value = 490000
$ws.Range("O574").Value = 510000
$ws.Range("P574").Value = 500000
$ws.Range("Q574").Value = "`$/bins (450 kilos)"
$ws.Range("R574").Value = "Región de O'Higgins"
$ws.Range("S574").Value = 1111
$ws.Range("T574").Value = 450
